$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ A=96997637; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540318.4948027555; R=6945610.198437067 }
  3  = @{ A=96997111; B=90669;  D='NT'; E=2059;   F='Skrovlig taggsvamp';   G='Hydnellum scabrosum';      H='(Fr.) E.Larss., K.H.Larss. & Kõljalg';   Q=540285.420278896;  R=6945753.778063174 }
  4  = @{ A=96997190; B=77259;  D='NT'; E=228912; F='Mörk kolflarnlav';     G='Carbonicola myrmecina';    H='(Ach.) Bendiksby & Timdal';               Q=540320.7145815691; R=6945654.846978122 }
  5  = @{ A=96997345; B=77258;  D='NT'; E=6446;   F='Kolflarnlav';          G='Carbonicola anthracophila'; H='(Nyl.) Bendiksby & Timdal';              Q=540324.0749594209; R=6945719.289026539 }
  6  = @{ A=96997189; B=77259;  D='NT'; E=228912; F='Mörk kolflarnlav';     G='Carbonicola myrmecina';    H='(Ach.) Bendiksby & Timdal';               Q=540324.0749594209; R=6945719.289026539 }
  7  = @{ A=96997467; B=77506;  D='NT'; E=6425;   F='Garnlav';              G='Alectoria sarmentosa';     H='(Ach.) Ach.';                             Q=540321.0604783226; R=6945550.886808772 }
  8  = @{ A=96997638; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540321.0604783226; R=6945550.886808772 }
  9  = @{ A=96997344; B=77258;  D='NT'; E=6446;   F='Kolflarnlav';          G='Carbonicola anthracophila'; H='(Nyl.) Bendiksby & Timdal';              Q=540239.4233446738; R=6945748.158169477 }
  10 = @{ A=96997470; B=77506;  D='NT'; E=6425;   F='Garnlav';              G='Alectoria sarmentosa';     H='(Ach.) Ach.';                             Q=540339.7901433307; R=6945337.195879548 }
  11 = @{ A=96997471; B=77506;  D='NT'; E=6425;   F='Garnlav';              G='Alectoria sarmentosa';     H='(Ach.) Ach.';                             Q=540403.0780395848; R=6945210.07007421 }
  12 = @{ A=96997642; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540357.5588053473; R=6945353.514534771 }
  13 = @{ A=96997641; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540387.8320511305; R=6945402.189333159 }
  14 = @{ A=96997643; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540338.9642116233; R=6945329.36484168 }
  15 = @{ A=96997135; B=90665;  D='LC'; E=4366;   F='Skarp dropptaggsvamp'; G='Hydnellum peckii';         H='Banker';                                  Q=540346.507806809;  R=6945277.469942174 }
  16 = @{ A=96997647; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540376.8706205683; R=6945167.883105326 }
  17 = @{ A=96997469; B=77506;  D='NT'; E=6425;   F='Garnlav';              G='Alectoria sarmentosa';     H='(Ach.) Ach.';                             Q=540389.3095380173; R=6945432.110323969 }
  18 = @{ A=96997646; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540384.6907192031; R=6945168.898715561 }
  19 = @{ A=96997644; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540341.1096913856; R=6945304.547811897 }
  20 = @{ A=96997645; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540342.1656137862; R=6945293.519209918 }
  21 = @{ A=96997640; B=90653;  D='LC'; E=4364;   F='Dropptaggsvamp';       G='Hydnellum ferrugineum';    H='(Fr.:Fr.) P. Karst.';                    Q=540389.3095380173; R=6945432.110323969 }
}

foreach ($rowNum in $data.Keys) {
  $rec = $data[$rowNum]
  $ws.Range("A$rowNum").Value = $rec.A
  $ws.Range("B$rowNum").Value = $rec.B
  $ws.Range("D$rowNum").Value = $rec.D
  $ws.Range("E$rowNum").Value = $rec.E
  $ws.Range("F$rowNum").Value = $rec.F
  $ws.Range("G$rowNum").Value = $rec.G
  $ws.Range("H$rowNum").Value = $rec.H
  $ws.Range("Q$rowNum").Value = $rec.Q
  $ws.Range("R$rowNum").Value = $rec.R
}
